# "Fruta / hortaliza, semanal" — weekly refresh: a new daily record is
# inserted at the top of the data (row 6) and all the previously existing
# records shift down by one row. The dimension grows from A1:T17 to A1:T18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6 (pushes rows 6-17 down to 7-18).
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the latest weekly record. All the
# "dimension" columns (market, region, product taxonomy, unit, origin,
# kg/unit) are identical to the rest of the series, only the date and the
# traded volume / derived price figures differ.
$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(6, 3).Value = "Coquimbo"

$ws.Cells.Item(6, 4).Value = 44413
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat

$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100108
$ws.Cells.Item(6, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(6, 9).Value = 100108007
$ws.Cells.Item(6, 10).Value = "Coco"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 45
$ws.Cells.Item(6, 14).Value = 20000
$ws.Cells.Item(6, 15).Value = 20000
$ws.Cells.Item(6, 16).Value = 20000
$ws.Cells.Item(6, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(6, 18).Value = "Perú"
$ws.Cells.Item(6, 19).Value = 1000
$ws.Cells.Item(6, 20).Value = 20
